# Enable return into game from any card, even if start outside of game
#
# This script:
#  1. Updates existing "count" (column B) values for a handful of cities and
#     records their original value in a new column C ("original").
#  2. Appends 26 new city rows (19-44) with their own counts.
#  3. Adds an AutoFilter over the full table plus the hidden
#     _xlnm._FilterDatabase defined name that Excel creates for it.
#  4. Tidies up cosmetic bits: column A width, and the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column C header + updates to existing rows that now carry an
#    "original" value alongside the new count.
# ---------------------------------------------------------------------------

$ws.Cells.Item(6, 2).Value = 2     # Lagos        B6  3 -> 2
$ws.Cells.Item(6, 3).Value = 3     # Lagos        C6  (original)

$ws.Cells.Item(7, 2).Value = 1     # Sao Paolo    B7  3 -> 1
$ws.Cells.Item(7, 3).Value = 2     # Sao Paolo    C7  (original)

$ws.Cells.Item(16, 2).Value = 0    # Mexico       B16 1 -> 0
$ws.Cells.Item(16, 3).Value = 1    # Mexico       C16 (original)

$ws.Cells.Item(17, 2).Value = 0    # Buenos Aires B17 2 -> 0
$ws.Cells.Item(17, 3).Value = 2    # Buenos Aires C17 (original)

$ws.Cells.Item(18, 2).Value = 0    # Lima         B18 1 -> 0

# ---------------------------------------------------------------------------
# 2. Append new rows 19-35. Insert them by copying the formatting of the
#    last existing data row (18) so column A keeps the same cell style,
#    then fill in the values.
# ---------------------------------------------------------------------------

for ($i = 19; $i -le 35; $i++) {
    $ws.Rows("18:18").Copy()
    $ws.Rows("$i`:$i").Insert(-4121)
}

# City names are assigned in this particular order so that the shared
# string table ends up with the same ordering as the source workbook.
$ws.Cells.Item(19, 1).Value = "Paris"
$ws.Cells.Item(21, 1).Value = "Moscou"
$ws.Cells.Item(22, 1).Value = "Saint-Petersbourg"
$ws.Cells.Item(23, 1).Value = "Francfort"
$ws.Cells.Item(24, 1).Value = "Bogota"
$ws.Cells.Item(20, 1).Value = "Santiago"
$ws.Cells.Item(25, 1).Value = "Kinshasa"
$ws.Cells.Item(26, 1).Value = "Antananarivo"
$ws.Cells.Item(27, 1).Value = "Dar es Salam"
$ws.Cells.Item(28, 1).Value = "Khartoum"
$ws.Cells.Item(29, 1).Value = "Johannesburg"
$ws.Cells.Item(30, 1).Value = "Bagdad"
$ws.Cells.Item(31, 1).Value = "Teheran"
$ws.Cells.Item(32, 1).Value = "Riyad"
$ws.Cells.Item(33, 1).Value = "Dehli"
$ws.Cells.Item(34, 1).Value = "Calcutta"
$ws.Cells.Item(35, 1).Value = "Navi Mumbai"

# Counts (column B) for rows 19-35.
$ws.Cells.Item(19, 2).Value = 2
$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(21, 2).Value = 1
$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(23, 2).Value = 2
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(25, 2).Value = 1
$ws.Cells.Item(26, 2).Value = 2
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(30, 2).Value = 2
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(34, 2).Value = 1
$ws.Cells.Item(35, 2).Value = 2

# Header for the new column (not bold, like the rest of the data cells).
# Assigned here so the shared string table keeps the same ordering as the
# source workbook.
$ws.Cells.Item(1, 3).Value = "original"

# ---------------------------------------------------------------------------
# 3. AutoFilter over the table (A1:B35), applied while that is still the
#    full extent of the data, plus the hidden worksheet-scoped
#    _xlnm._FilterDatabase defined name that Excel records for it.
# ---------------------------------------------------------------------------

$ws.Range("A1:B35").AutoFilter(1) | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$B`$35")
$filterName.Visible = $false

# ---------------------------------------------------------------------------
# 4. Append the remaining rows 36-44.
# ---------------------------------------------------------------------------

for ($i = 36; $i -le 44; $i++) {
    $ws.Rows("18:18").Copy()
    $ws.Rows("$i`:$i").Insert(-4121)
}

$ws.Cells.Item(36, 1).Value = "Jakarta"
$ws.Cells.Item(37, 1).Value = "Bankok"
$ws.Cells.Item(38, 1).Value = "Ho-Chi-Minh-Ville"
$ws.Cells.Item(39, 1).Value = "Manille"
$ws.Cells.Item(40, 1).Value = "Seoul"
$ws.Cells.Item(41, 1).Value = "Tokyo"
$ws.Cells.Item(42, 1).Value = "Osaka"
$ws.Cells.Item(43, 1).Value = "Shanghai"
$ws.Cells.Item(44, 1).Value = "Hong Kong"

$ws.Cells.Item(36, 2).Value = 1
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(38, 2).Value = 0
$ws.Cells.Item(39, 2).Value = 0
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(43, 2).Value = 0
$ws.Cells.Item(44, 2).Value = 0

# ---------------------------------------------------------------------------
# 5. Cosmetics: widen column A to fit the longer city names, and leave the
#    selection on B1.
# ---------------------------------------------------------------------------

$ws.Columns("A:A").ColumnWidth = 14.7
$ws.Range("B1").Select() | Out-Null
